# Generate Report for Handoff
# A new handoff cycle was generated for the source file, producing a new
# GUID-based file name and new xliff checksums. The new handback has not
# completed yet for either target locale, so the "Latest Target File" /
# "Latest Handback File" columns (and their hyperlinks) are cleared and the
# "Latest Handback DateTime" reverts to the unset default.

$wb = $excel.ActiveWorkbook

$oldGuid = "ea948818-6ad9-446f-b6b8-e4bae19996e3"
$newGuid = "7b5ff3fd-bf33-4231-b570-582bdfdcae07"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "$newGuid.md"
$overview.Range("B2").Value = "e2e\$newGuid.md"
$overview.Range("G2").Value = "2016-08-27 00:59:33"
$overview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "$newGuid.md"
$zhcn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$zhcn.Range("G2").Value = "$newGuid.aed12fa832315da62399b9d1eddc68662799a56b.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-27 00:59:28"

# Handback hasn't happened for this generation yet: drop the stale
# target/handback file hyperlink and blank out the dependent columns.
$zhcn.Hyperlinks.Item(2).Delete()
$zhcn.Range("I2").Value = ""
$zhcn.Range("I2").Style = "Normal"
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "$newGuid.md"
$dede.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$dede.Range("G2").Value = "$newGuid.aed12fa832315da62399b9d1eddc68662799a56b.de-de.xlf"
$dede.Range("H2").Value = "2016-08-27 00:59:33"

$dede.Hyperlinks.Item(2).Delete()
$dede.Range("I2").Value = ""
$dede.Range("I2").Style = "Normal"
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = "0001-01-01 00:00:00"
